$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Brand" column (C) was incorrectly populated with the Generic placeholder.
# Fix it so it uses the correct Brand placeholder for both data rows.
$ws.Range("C2").Value = "{d.records[i].brand}"
$ws.Range("C3").Value = "{d.records[i+1].brand}"

$ws.Range("C3").Select()
